$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rina")

# Insert a blank column before K. K:K sits inside the uniform-width C:K
# band, so the new column merges into that band's width exactly, and the
# old K ("WiLAN") column's data/styles shift one column right, to L.
$ws.Columns("K:K").Insert()

# Put WiLAN's data back onto K (styles already followed the shift).
$ws.Cells.Item(2, 11).Value = "WiLAN"
$ws.Cells.Item(3, 11).Value = "2/0/3"
$ws.Cells.Item(4, 11).Value = 5
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(6, 11).Value = 4
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(8, 11).Formula = "=K5+K6+K7"
$ws.Cells.Item(9, 11).Value = 13
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(11, 11).Formula = "=K9+K10"
$ws.Cells.Item(12, 11).Formula = "=K4+K8"
$ws.Cells.Item(13, 11).Formula = "=K8+K11"
$ws.Cells.Item(14, 11).Formula = "=K4+K11"
$ws.Cells.Item(15, 11).Formula = "=K4+K8+K11"
$ws.Cells.Item(16, 11).Value = 8
$ws.Cells.Item(17, 11).Value = 4
$ws.Cells.Item(18, 11).Formula = "=K16+K17"
$ws.Cells.Item(19, 11).Formula = "=K15+K18"

# New "OS" column (L): header, sub-header, and data.
$ws.Cells.Item(2, 12).Value = "OS"
$ws.Cells.Item(3, 12).Value = "2/0/1"
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(6, 12).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(8, 12).Formula = "=L5+L6+L7"
$ws.Cells.Item(9, 12).Value = 4
$ws.Cells.Item(10, 12).Value = 3
$ws.Cells.Item(11, 12).Formula = "=L9+L10"
$ws.Cells.Item(12, 12).Formula = "=L4+L8"
$ws.Cells.Item(13, 12).Formula = "=L8+L11"
$ws.Cells.Item(14, 12).Formula = "=L4+L11"
$ws.Cells.Item(15, 12).Formula = "=L4+L8+L11"
$ws.Cells.Item(16, 12).Value = 4
$ws.Cells.Item(17, 12).Value = 2
$ws.Cells.Item(18, 12).Formula = "=L16+L17"
$ws.Cells.Item(19, 12).Formula = "=L15+L18"

# Extend the print area to cover the new last column (M, the right margin).
$ws.PageSetup.PrintArea = '$A$1:$M$20'

# Move the selection to L18, matching the saved view state.
$ws.Range("L18").Select()
